$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grade")
$ws.Activate()

# C18: 0 -> 3
$ws.Range("C18").Value = 3

# D18: remove the "Not used" label entirely
$ws.Range("D18").ClearContents()

# Scroll the view down to row 13 and zoom to 150%, then select E18
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 150
$ws.Range("E18").Select()

Write-Output "done"
